$d = $word.ActiveDocument

# Update title date paragraph (2023-07-22 Saturday -> 2023-07-23 Sunday)
$d.Paragraphs(1).Range.Text = "2023-07-23 Sunday"

# Update each answer cell in the multiplication table, in row-major order
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "97×46=4462"
$t.Cell(1,2).Range.Text = "15×64=960"
$t.Cell(1,3).Range.Text = "10×94=940"
$t.Cell(1,4).Range.Text = "32×49=1568"
$t.Cell(1,5).Range.Text = "38×75=2850"
$t.Cell(2,1).Range.Text = "43×89=3827"
$t.Cell(2,2).Range.Text = "70×14=980"
$t.Cell(2,3).Range.Text = "38×47=1786"
$t.Cell(2,4).Range.Text = "38×46=1748"
$t.Cell(2,5).Range.Text = "33×74=2442"
$t.Cell(3,1).Range.Text = "69×43=2967"
$t.Cell(3,2).Range.Text = "98×18=1764"
$t.Cell(3,3).Range.Text = "98×32=3136"
$t.Cell(3,4).Range.Text = "42×18=756"
$t.Cell(3,5).Range.Text = "23×18=414"
$t.Cell(4,1).Range.Text = "32×54=1728"
$t.Cell(4,2).Range.Text = "80×16=1280"
$t.Cell(4,3).Range.Text = "28×55=1540"
$t.Cell(4,4).Range.Text = "86×60=5160"
$t.Cell(4,5).Range.Text = "88×32=2816"
$t.Cell(5,1).Range.Text = "40×43=1720"
$t.Cell(5,2).Range.Text = "24×36=864"
$t.Cell(5,3).Range.Text = "33×82=2706"
$t.Cell(5,4).Range.Text = "97×13=1261"
$t.Cell(5,5).Range.Text = "28×25=700"
$t.Cell(6,1).Range.Text = "29×37=1073"
$t.Cell(6,2).Range.Text = "22×23=506"
$t.Cell(6,3).Range.Text = "62×13=806"
$t.Cell(6,4).Range.Text = "59×67=3953"
$t.Cell(6,5).Range.Text = "46×63=2898"
$t.Cell(7,1).Range.Text = "79×60=4740"
$t.Cell(7,2).Range.Text = "88×15=1320"
$t.Cell(7,3).Range.Text = "68×40=2720"
$t.Cell(7,4).Range.Text = "52×46=2392"
$t.Cell(7,5).Range.Text = "48×30=1440"
$t.Cell(8,1).Range.Text = "39×39=1521"
$t.Cell(8,2).Range.Text = "60×91=5460"
$t.Cell(8,3).Range.Text = "73×95=6935"
$t.Cell(8,4).Range.Text = "35×65=2275"
$t.Cell(8,5).Range.Text = "78×20=1560"
$t.Cell(9,1).Range.Text = "94×49=4606"
$t.Cell(9,2).Range.Text = "34×63=2142"
$t.Cell(9,3).Range.Text = "11×73=803"
$t.Cell(9,4).Range.Text = "71×59=4189"
$t.Cell(9,5).Range.Text = "57×23=1311"
$t.Cell(10,1).Range.Text = "43×95=4085"
$t.Cell(10,2).Range.Text = "26×80=2080"
$t.Cell(10,3).Range.Text = "52×34=1768"
$t.Cell(10,4).Range.Text = "57×46=2622"
$t.Cell(10,5).Range.Text = "37×33=1221"
$t.Cell(11,1).Range.Text = "76×67=5092"
$t.Cell(11,2).Range.Text = "23×85=1955"
$t.Cell(11,3).Range.Text = "47×12=564"
$t.Cell(11,4).Range.Text = "21×31=651"
$t.Cell(11,5).Range.Text = "12×61=732"
$t.Cell(12,1).Range.Text = "16×12=192"
$t.Cell(12,2).Range.Text = "79×46=3634"
$t.Cell(12,3).Range.Text = "10×15=150"
$t.Cell(12,4).Range.Text = "16×21=336"
$t.Cell(12,5).Range.Text = "13×85=1105"
$t.Cell(13,1).Range.Text = "44×96=4224"
$t.Cell(13,2).Range.Text = "34×73=2482"
$t.Cell(13,3).Range.Text = "57×62=3534"
$t.Cell(13,4).Range.Text = "20×92=1840"
$t.Cell(13,5).Range.Text = "34×63=2142"
$t.Cell(14,1).Range.Text = "83×98=8134"
$t.Cell(14,2).Range.Text = "30×86=2580"
$t.Cell(14,3).Range.Text = "92×12=1104"
$t.Cell(14,4).Range.Text = "41×79=3239"
$t.Cell(14,5).Range.Text = "67×97=6499"
$t.Cell(15,1).Range.Text = "100×60=6000"
$t.Cell(15,2).Range.Text = "62×69=4278"
$t.Cell(15,3).Range.Text = "82×99=8118"
$t.Cell(15,4).Range.Text = "59×51=3009"
$t.Cell(15,5).Range.Text = "82×75=6150"
$t.Cell(16,1).Range.Text = "51×35=1785"
$t.Cell(16,2).Range.Text = "66×75=4950"
$t.Cell(16,3).Range.Text = "24×81=1944"
$t.Cell(16,4).Range.Text = "100×93=9300"
$t.Cell(16,5).Range.Text = "58×84=4872"
$t.Cell(17,1).Range.Text = "30×56=1680"
$t.Cell(17,2).Range.Text = "48×47=2256"
$t.Cell(17,3).Range.Text = "15×11=165"
$t.Cell(17,4).Range.Text = "80×42=3360"
$t.Cell(17,5).Range.Text = "27×25=675"
$t.Cell(18,1).Range.Text = "46×14=644"
$t.Cell(18,2).Range.Text = "56×82=4592"
$t.Cell(18,3).Range.Text = "93×11=1023"
$t.Cell(18,4).Range.Text = "92×21=1932"
$t.Cell(18,5).Range.Text = "63×86=5418"
$t.Cell(19,1).Range.Text = "62×35=2170"
$t.Cell(19,2).Range.Text = "66×100=6600"
$t.Cell(19,3).Range.Text = "56×79=4424"
$t.Cell(19,4).Range.Text = "93×90=8370"
$t.Cell(19,5).Range.Text = "14×78=1092"
$t.Cell(20,1).Range.Text = "87×41=3567"
$t.Cell(20,2).Range.Text = "93×37=3441"
$t.Cell(20,3).Range.Text = "75×37=2775"
$t.Cell(20,4).Range.Text = "96×11=1056"
$t.Cell(20,5).Range.Text = "21×26=546"
